$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNoDate($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextNoDate $ws.Range("C2") "01/07/2014"
$ws.Range("D2").Value = 79.04766161593587

Set-TextNoDate $ws.Range("C3") "01/07/2016"
$ws.Range("D3").Value = 79.67335137194378

Set-TextNoDate $ws.Range("C4") "01/07/2018"
$ws.Range("D4").Value = 80.16927727925945

Set-TextNoDate $ws.Range("C5") "01/07/2020"
$ws.Range("D5").Value = 80.596605163257

Set-TextNoDate $ws.Range("C6") "01/07/2022"
$ws.Range("D6").Value = 80.9366717682631

Set-TextNoDate $ws.Range("C7") "01/07/2024"
$ws.Range("D7").Value = 81.2474962126619

Set-TextNoDate $ws.Range("C8") "01/07/2014"
$ws.Range("D8").Value = 49.27516798522028

Set-TextNoDate $ws.Range("C9") "01/07/2016"
$ws.Range("D9").Value = 49.80707803815459

Set-TextNoDate $ws.Range("C10") "01/07/2018"
$ws.Range("D10").Value = 50.76588115985216

Set-TextNoDate $ws.Range("C11") "01/07/2020"
$ws.Range("D11").Value = 46.40541129024624

Set-TextNoDate $ws.Range("C12") "01/07/2022"
$ws.Range("D12").Value = 50.73397663219978

Set-TextNoDate $ws.Range("C13") "01/07/2024"
$ws.Range("D13").Value = 50.66514405698735

Set-TextNoDate $ws.Range("C14") "01/07/2014"
$ws.Range("D14").Value = 45.89212195255193

Set-TextNoDate $ws.Range("C15") "01/07/2016"
$ws.Range("D15").Value = 43.86947730358591

Set-TextNoDate $ws.Range("C16") "01/07/2018"
$ws.Range("D16").Value = 44.66478580800823

Set-TextNoDate $ws.Range("C17") "01/07/2020"
$ws.Range("D17").Value = 39.49550794747753

Set-TextNoDate $ws.Range("C18") "01/07/2022"
$ws.Range("D18").Value = 46.31985143155773

Set-TextNoDate $ws.Range("C19") "01/07/2024"
$ws.Range("D19").Value = 47.44141713212169

Set-TextNoDate $ws.Range("C20") "01/07/2014"
$ws.Range("D20").Value = 3.383046032668345

Set-TextNoDate $ws.Range("C21") "01/07/2016"
$ws.Range("D21").Value = 5.937112324538696

Set-TextNoDate $ws.Range("C22") "01/07/2018"
$ws.Range("D22").Value = 6.101095351843931

Set-TextNoDate $ws.Range("C23") "01/07/2020"
$ws.Range("D23").Value = 6.909903342768694

Set-TextNoDate $ws.Range("C24") "01/07/2022"
$ws.Range("D24").Value = 4.414125200642054

Set-TextNoDate $ws.Range("C25") "01/07/2024"
$ws.Range("D25").Value = 3.223726924865659

Set-TextNoDate $ws.Range("C26") "01/07/2014"
$ws.Range("D26").Value = 29.7724936307156

Set-TextNoDate $ws.Range("C27") "01/07/2016"
$ws.Range("D27").Value = 29.86676174381917

Set-TextNoDate $ws.Range("C28") "01/07/2018"
$ws.Range("D28").Value = 29.40339611940729

Set-TextNoDate $ws.Range("C29") "01/07/2020"
$ws.Range("D29").Value = 34.19166721890355

Set-TextNoDate $ws.Range("C30") "01/07/2022"
$ws.Range("D30").Value = 30.20316174549255

Set-TextNoDate $ws.Range("C31") "01/07/2024"
$ws.Range("D31").Value = 30.58235215567456

Set-TextNoDate $ws.Range("C32") "01/07/2014"
$ws.Range("D32").Value = 76.77104954426495

Set-TextNoDate $ws.Range("C33") "01/07/2016"
$ws.Range("D33").Value = 77.96395431834404

Set-TextNoDate $ws.Range("C34") "01/07/2018"
$ws.Range("D34").Value = 78.62626405237985

Set-TextNoDate $ws.Range("C35") "01/07/2020"

Set-TextNoDate $ws.Range("C36") "01/07/2022"
$ws.Range("D36").Value = 79.771414738735

Set-TextNoDate $ws.Range("C37") "01/07/2024"
$ws.Range("D37").Value = 80.25357662982952

Set-TextNoDate $ws.Range("C38") "01/07/2014"
$ws.Range("D38").Value = 44.51764281202058

Set-TextNoDate $ws.Range("C39") "01/07/2016"
$ws.Range("D39").Value = 43.76695217701641

Set-TextNoDate $ws.Range("C40") "01/07/2018"
$ws.Range("D40").Value = 44.13991493567231

Set-TextNoDate $ws.Range("C41") "01/07/2020"

Set-TextNoDate $ws.Range("C42") "01/07/2022"
$ws.Range("D42").Value = 43.96375834284331

Set-TextNoDate $ws.Range("C43") "01/07/2024"
$ws.Range("D43").Value = 43.92015762871584

Set-TextNoDate $ws.Range("C44") "01/07/2014"
$ws.Range("D44").Value = 40.65517552567458

Set-TextNoDate $ws.Range("C45") "01/07/2016"
$ws.Range("D45").Value = 37.52498215560314

Set-TextNoDate $ws.Range("C46") "01/07/2018"
$ws.Range("D46").Value = 37.70185129626034

Set-TextNoDate $ws.Range("C47") "01/07/2020"

Set-TextNoDate $ws.Range("C48") "01/07/2022"
$ws.Range("D48").Value = 38.70041843898053

Set-TextNoDate $ws.Range("C49") "01/07/2024"
$ws.Range("D49").Value = 40.11136811445215

Set-TextNoDate $ws.Range("C50") "01/07/2014"
$ws.Range("D50").Value = 3.862467286345998

Set-TextNoDate $ws.Range("C51") "01/07/2016"
$ws.Range("D51").Value = 6.241970021413276

Set-TextNoDate $ws.Range("C52") "01/07/2018"
$ws.Range("D52").Value = 6.436298819335369

Set-TextNoDate $ws.Range("C53") "01/07/2020"

Set-TextNoDate $ws.Range("C54") "01/07/2022"
$ws.Range("D54").Value = 5.26333990386278

Set-TextNoDate $ws.Range("C55") "01/07/2024"
$ws.Range("D55").Value = 3.808789514263686

Set-TextNoDate $ws.Range("C56") "01/07/2014"
$ws.Range("D56").Value = 32.25340673224439

Set-TextNoDate $ws.Range("C57") "01/07/2016"
$ws.Range("D57").Value = 34.19700214132762

Set-TextNoDate $ws.Range("C58") "01/07/2018"
$ws.Range("D58").Value = 34.48811393678415

Set-TextNoDate $ws.Range("C59") "01/07/2020"

Set-TextNoDate $ws.Range("C60") "01/07/2022"
$ws.Range("D60").Value = 35.80765639589169

Set-TextNoDate $ws.Range("C61") "01/07/2024"
$ws.Range("D61").Value = 36.33341900111368

Set-TextNoDate $ws.Range("C62") "01/07/2014"
$ws.Range("D62").Value = 76.55894401456531

Set-TextNoDate $ws.Range("C63") "01/07/2016"
$ws.Range("D63").Value = 78.83928571428571

Set-TextNoDate $ws.Range("C64") "01/07/2018"
$ws.Range("D64").Value = 77.77290661990355

Set-TextNoDate $ws.Range("C65") "01/07/2020"

Set-TextNoDate $ws.Range("C66") "01/07/2022"
$ws.Range("D66").Value = 79.11016949152543

Set-TextNoDate $ws.Range("C67") "01/07/2024"
$ws.Range("D67").Value = 79.96661101836395

Set-TextNoDate $ws.Range("C68") "01/07/2014"
$ws.Range("D68").Value = 46.15384615384615

Set-TextNoDate $ws.Range("C69") "01/07/2016"
$ws.Range("D69").Value = 45.80357142857142

Set-TextNoDate $ws.Range("C70") "01/07/2018"
$ws.Range("D70").Value = 45.46251644015783

Set-TextNoDate $ws.Range("C71") "01/07/2020"

Set-TextNoDate $ws.Range("C72") "01/07/2022"
$ws.Range("D72").Value = 45.76271186440678

Set-TextNoDate $ws.Range("C73") "01/07/2024"
$ws.Range("D73").Value = 46.78631051752922

Set-TextNoDate $ws.Range("C74") "01/07/2014"
$ws.Range("D74").Value = 41.92080109239873

Set-TextNoDate $ws.Range("C75") "01/07/2016"
$ws.Range("D75").Value = 39.15178571428572

Set-TextNoDate $ws.Range("C76") "01/07/2018"
$ws.Range("D76").Value = 37.48355984217449

Set-TextNoDate $ws.Range("C77") "01/07/2020"

Set-TextNoDate $ws.Range("C78") "01/07/2022"
$ws.Range("D78").Value = 40.21186440677966

Set-TextNoDate $ws.Range("C79") "01/07/2024"
$ws.Range("D79").Value = 42.86310517529216

Set-TextNoDate $ws.Range("C80") "01/07/2014"
$ws.Range("D80").Value = 4.233045061447428

Set-TextNoDate $ws.Range("C81") "01/07/2016"
$ws.Range("D81").Value = 6.607142857142858

Set-TextNoDate $ws.Range("C82") "01/07/2018"
$ws.Range("D82").Value = 8.022797018851382

Set-TextNoDate $ws.Range("C83") "01/07/2020"

Set-TextNoDate $ws.Range("C84") "01/07/2022"
$ws.Range("D84").Value = 5.550847457627119

Set-TextNoDate $ws.Range("C85") "01/07/2024"
$ws.Range("D85").Value = 3.923205342237062

Set-TextNoDate $ws.Range("C86") "01/07/2014"
$ws.Range("D86").Value = 30.40509786071916

Set-TextNoDate $ws.Range("C87") "01/07/2016"
$ws.Range("D87").Value = 33.03571428571428

Set-TextNoDate $ws.Range("C88") "01/07/2018"
$ws.Range("D88").Value = 32.26654975887768

Set-TextNoDate $ws.Range("C89") "01/07/2020"

Set-TextNoDate $ws.Range("C90") "01/07/2022"
$ws.Range("D90").Value = 33.34745762711864

Set-TextNoDate $ws.Range("C91") "01/07/2024"
$ws.Range("D91").Value = 33.13856427378965
